# Update the emission cap year values: 2020 -> 2022, 2025 -> 2026
# Applies to both the "emission_cap" sheet and the "emission_cap_old" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("emission_cap")
$ws1.Range("A2").Value = 2022
$ws1.Range("A3").Value = 2026

$ws2 = $wb.Worksheets.Item("emission_cap_old")
$ws2.Range("A2").Value = 2022
$ws2.Range("A3").Value = 2026

# Restore the active sheet / selection to mirror the recorded cursor positions.
$ws1.Activate()
$ws1.Range("E29").Select()

$ws2.Activate()
$ws2.Range("A4").Select()

$ws1.Activate()
